$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $c = $ws.Range($cell)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "67.904.33"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "3.740.63"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue "D5" "593.72"
$ws.Range("E5").Value = "  -1.13%  "
Set-TextValue "D6" "166.63"
$ws.Range("E6").Value = "  -1.93%  "
$ws.Range("D7").Value = "3.741.19"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -1.98%  "
$ws.Range("E10").Value = "  -3.92%  "
Set-TextValue "D11" "6.52"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("E12").Value = "  -2.37%  "
Set-TextValue "D13" "0.0000264"
$ws.Range("E13").Value = "  -5.04%  "
Set-TextValue "D14" "36.42"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "4.364.47"
$ws.Range("E15").Value = "  -2.22%  "
$ws.Range("D16").Value = "3.743.87"
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "67.886.45"
$ws.Range("E17").Value = "  -0.17%  "
Set-TextValue "D18" "18.27"
$ws.Range("E18").Value = "  -2.34%  "
Set-TextValue "D19" "7.05"
$ws.Range("E19").Value = "  -5.53%  "
$ws.Range("E20").Value = "  -0.45%  "
Set-TextValue "D21" "10.85"
$ws.Range("E21").Value = "  -0.52%  "
Set-TextValue "D22" "467.45"
$ws.Range("E22").Value = "  -0.54%  "
Set-TextValue "D23" "0.702"
$ws.Range("E23").Value = "  -5.40%  "
Set-TextValue "D24" "82.97"
$ws.Range("E25").Value = "  -3.88%  "
$ws.Range("E26").Value = "  -11.30%  "
$ws.Range("E27").Value = "  -1.63%  "
Set-TextValue "D28" "10.17"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "3.885.55"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("E31").Value = "  -4.96%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "2.24"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D33" "7.35"
$ws.Range("E33").Value = "  -5.33%  "
Set-TextValue "D34" "29.79"
$ws.Range("E34").Value = "  -3.29%  "
Set-TextValue "D35" "1.00"
Set-TextValue "D36" "9.09"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("D37").Value = "3.691.29"
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("E38").Value = "  -4.54%  "
Set-TextValue "D39" "3.43"
$ws.Range("E39").Value = "  -11.86%  "
$ws.Range("E40").Value = "  -1.24%  "
Set-TextValue "D41" "0.995"
$ws.Range("E41").Value = "  -1.60%  "
Set-TextValue "D42" "5.76"
$ws.Range("E42").Value = "  -4.07%  "
Set-TextValue "D43" "0.999"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -4.01%  "
Set-TextValue "D46" "8.58"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("E47").Value = "  -3.15%  "
Set-TextValue "D48" "45.41"
$ws.Range("E48").Value = "  -2.53%  "
Set-TextValue "D49" "391.55"
$ws.Range("E49").Value = "  -5.13%  "
Set-TextValue "D50" "143.57"
$ws.Range("E50").Value = "  +0.56%  "
Set-TextValue "D51" "25.50"
$ws.Range("E51").Value = "  +0.27%  "
